$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture B4's current format (the "last task" marker style) before we change it, so we can
# reapply it to the new row 14 below.
$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats

# Row 2: rename the "Lesson wizard" task to its stage #1, and highlight it as in-progress
# (copy the highlight format from B1, which already carries the "in progress" style).
$ws.Range("B2").Value2 = 'Create wizard for "Lesson" page stage #1'
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 ("Change Karma test browser to PhantomJS") becomes highlighted too, same style as B1/B2.
$ws.Range("B1").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats

# New row 14: the "Lesson wizard" stage #2 task (format was already applied above).
$ws.Range("B14").Value2 = 'Create wizard for "Lesson" page stage #2'

# Tab color switch from explicit RGB green to the equivalent theme color (theme 6, tint ~0.6)
$ws.Tab.ThemeColor = 6
$ws.Tab.TintAndShade = 0.59999389629810485

# Selection moves to C13
$ws.Range("C13").Select()
